# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Thu Jun 13 10:52:16 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: Coin (B), Link (C), Price (D) and Volume(1h) (E) columns.
# Column D holds price strings that can look numeric (e.g. "1.00", "151.30"),
# so those cells are forced to Text format first to keep the exact string
# (matching the original inline-string cell content) instead of Excel
# silently re-parsing them as numbers and dropping trailing zeros.
$updates = @(
    @{ Cell = "D2"; Value = "67.629.18"; AsText = $true }
    @{ Cell = "E2"; Value = "  -0.36%  "; AsText = $false }
    @{ Cell = "D3"; Value = "3.498.52"; AsText = $true }
    @{ Cell = "E3"; Value = "  -1.14%  "; AsText = $false }
    @{ Cell = "E4"; Value = "  -0.15%  "; AsText = $false }
    @{ Cell = "D5"; Value = "606.95"; AsText = $true }
    @{ Cell = "E5"; Value = "  -1.45%  "; AsText = $false }
    @{ Cell = "D6"; Value = "151.30"; AsText = $true }
    @{ Cell = "E6"; Value = "  -0.77%  "; AsText = $false }
    @{ Cell = "D7"; Value = "3.494.71"; AsText = $true }
    @{ Cell = "E7"; Value = "  -1.19%  "; AsText = $false }
    @{ Cell = "D8"; Value = "1.00"; AsText = $true }
    @{ Cell = "E8"; Value = "  -0.06%  "; AsText = $false }
    @{ Cell = "D9"; Value = "0.487"; AsText = $true }
    @{ Cell = "E9"; Value = "  +0.61%  "; AsText = $false }
    @{ Cell = "E10"; Value = "  +2.76%  "; AsText = $false }
    @{ Cell = "D11"; Value = "7.58"; AsText = $true }
    @{ Cell = "E11"; Value = "  +6.89%  "; AsText = $false }
    @{ Cell = "E12"; Value = "  +1.02%  "; AsText = $false }
    @{ Cell = "E13"; Value = "  -2.15%  "; AsText = $false }
    @{ Cell = "D14"; Value = "32.04"; AsText = $true }
    @{ Cell = "E14"; Value = "  -0.31%  "; AsText = $false }
    @{ Cell = "D15"; Value = "4.086.14"; AsText = $true }
    @{ Cell = "E15"; Value = "  -1.30%  "; AsText = $false }
    @{ Cell = "D16"; Value = "67.606.90"; AsText = $true }
    @{ Cell = "E16"; Value = "  -0.03%  "; AsText = $false }
    @{ Cell = "D17"; Value = "3.487.72"; AsText = $true }
    @{ Cell = "E17"; Value = "  -1.53%  "; AsText = $false }
    @{ Cell = "E18"; Value = "  -0.19%  "; AsText = $false }
    @{ Cell = "D19"; Value = "6.49"; AsText = $true }
    @{ Cell = "E19"; Value = "  +1.28%  "; AsText = $false }
    @{ Cell = "D20"; Value = "15.42"; AsText = $true }
    @{ Cell = "D21"; Value = "9.92"; AsText = $true }
    @{ Cell = "E21"; Value = "  +2.04%  "; AsText = $false }
    @{ Cell = "D22"; Value = "446.67"; AsText = $true }
    @{ Cell = "E22"; Value = "  -0.28%  "; AsText = $false }
    @{ Cell = "E23"; Value = "  +0.14%  "; AsText = $false }
    @{ Cell = "D24"; Value = "79.41"; AsText = $true }
    @{ Cell = "E24"; Value = "  +2.31%  "; AsText = $false }
    @{ Cell = "D25"; Value = "3.634.19"; AsText = $true }
    @{ Cell = "E25"; Value = "  -1.31%  "; AsText = $false }
    @{ Cell = "E26"; Value = "  +0.03%  "; AsText = $false }
    @{ Cell = "E27"; Value = "  -4.10%  "; AsText = $false }
    @{ Cell = "E28"; Value = "  +0.23%  "; AsText = $false }
    @{ Cell = "D29"; Value = "9.94"; AsText = $true }
    @{ Cell = "E29"; Value = "  -3.24%  "; AsText = $false }
    @{ Cell = "D30"; Value = "2.51"; AsText = $true }
    @{ Cell = "E30"; Value = "  -1.26%  "; AsText = $false }
    @{ Cell = "D31"; Value = "1.65"; AsText = $true }
    @{ Cell = "E31"; Value = "  +2.67%  "; AsText = $false }
    @{ Cell = "E32"; Value = "  +1.50%  "; AsText = $false }
    @{ Cell = "D33"; Value = "0.997"; AsText = $true }
    @{ Cell = "E33"; Value = "  -0.18%  "; AsText = $false }
    @{ Cell = "D34"; Value = "25.63"; AsText = $true }
    @{ Cell = "E34"; Value = "  -1.35%  "; AsText = $false }
    @{ Cell = "D35"; Value = "6.14"; AsText = $true }
    @{ Cell = "E35"; Value = "  -1.29%  "; AsText = $false }
    @{ Cell = "E36"; Value = "  -0.17%  "; AsText = $false }
    @{ Cell = "D37"; Value = "3.490.49"; AsText = $true }
    @{ Cell = "E37"; Value = "  -1.06%  "; AsText = $false }
    @{ Cell = "D38"; Value = "8.02"; AsText = $true }
    @{ Cell = "E38"; Value = "  -0.53%  "; AsText = $false }
    @{ Cell = "E39"; Value = "  +0.02%  "; AsText = $false }
    @{ Cell = "D40"; Value = "2.32"; AsText = $true }
    @{ Cell = "E40"; Value = "  +5.59%  "; AsText = $false }
    @{ Cell = "B41"; Value = "Monero"; AsText = $false }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; AsText = $false }
    @{ Cell = "D41"; Value = "177.47"; AsText = $true }
    @{ Cell = "E41"; Value = "  +0.34%  "; AsText = $false }
    @{ Cell = "B42"; Value = "FirstDigitalUSD"; AsText = $false }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; AsText = $false }
    @{ Cell = "D42"; Value = "0.999"; AsText = $true }
    @{ Cell = "E42"; Value = "  -0.13%  "; AsText = $false }
    @{ Cell = "E43"; Value = "  +0.32%  "; AsText = $false }
    @{ Cell = "E44"; Value = "  +0.06%  "; AsText = $false }
    @{ Cell = "B45"; Value = "InjectiveProtocol"; AsText = $false }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; AsText = $false }
    @{ Cell = "D45"; Value = "30.46"; AsText = $true }
    @{ Cell = "E45"; Value = "  +6.73%  "; AsText = $false }
    @{ Cell = "B46"; Value = "Mantle"; AsText = $false }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; AsText = $false }
    @{ Cell = "D46"; Value = "0.892"; AsText = $true }
    @{ Cell = "E46"; Value = "  +0.68%  "; AsText = $false }
    @{ Cell = "D47"; Value = "46.51"; AsText = $true }
    @{ Cell = "E47"; Value = "  +2.51%  "; AsText = $false }
    @{ Cell = "E48"; Value = "  -0.13%  "; AsText = $false }
    @{ Cell = "E49"; Value = "  -5.04%  "; AsText = $false }
    @{ Cell = "E50"; Value = "  -0.24%  "; AsText = $false }
    @{ Cell = "E51"; Value = "  -0.75%  "; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

